$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "Angular.isNumber" -> "angular.isNumber", split into two runs
#    ("a" + "ngular.isNumber") the way Word splits a run when the
#    cursor sits in the middle of it. We use a transient bookmark to
#    force the structural split, then fix the casing of the first
#    character.
# -----------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Angular.isNumber", $true, $false, $false, $false, $false, $true, 1, $false)
$wordStart = $r.Start

$splitPoint = $d.Range($wordStart + 1, $wordStart + 1)
$d.Bookmarks.Add("_GoBack", $splitPoint)

$firstChar = $d.Range($wordStart, $wordStart + 1)
$firstChar.Text = "a"

# -----------------------------------------------------------------
# 2) "constant" -> "const" right before the closing curly quote in
#    "...using the keyword \u201cconstant\u201d." and leave the
#    document's _GoBack bookmark at the split point created by that
#    edit (Word always keeps _GoBack at the location of the most
#    recent edit, re-adding a bookmark with the same name simply
#    relocates it). The run is first split at the edit point (this is
#    also where the final bookmark belongs), then the trailing "ant"
#    is removed from the back half.
# -----------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("constant”. They are created using constant dependency", $true, $false, $false, $false, $false, $true, 1, $false)
$kwStart = $r2.Start

$finalBookmarkPoint = $d.Range($kwStart + 5, $kwStart + 5)
$d.Bookmarks.Add("_GoBack", $finalBookmarkPoint)

$suffix = $d.Range($kwStart + 5, $kwStart + 8)
$suffix.Text = ""
